$d = $word.ActiveDocument

function Get-ParagraphAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $pos -and $pos -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Remove the "Output her er veldig midlertidig..." paragraph together with
#    the empty paragraph that immediately follows it.
# ---------------------------------------------------------------------------
$find1 = $d.Content.Duplicate
$find1.Find.ClearFormatting()
$find1.Find.MatchCase = $true
$find1.Find.Text = "Output her er veldig midlertidig, og vil sannsynligvis til å forandres, da den går på vår gamle og utdaterte mal. Det vil allikevel gi en pekepinn på hva slags info som er viktig."
$found1 = $find1.Find.Execute()

if ($found1) {
    $targetPara = Get-ParagraphAt $d $find1.Start
    $nextPara = $targetPara.Next()
    $deleteRange = $d.Range($targetPara.Range.Start, $nextPara.Range.End)
    $deleteRange.Delete()
}

# ---------------------------------------------------------------------------
# 2. Rework the BaseX/"for opptelling av registreringer per år" sentence so
#    that the "Fordi Output er under omskriving..." aside is removed while
#    keeping the remaining runs intact (å sjekke / antallet... / the final
#    period untouched).
# ---------------------------------------------------------------------------
$tailStart = $d.Content.Duplicate
$tailStart.Find.ClearFormatting()
$tailStart.Find.MatchCase = $true
$tailStart.Find.Text = " for opptelling av registreringer per år. Fordi Output "
$foundTail = $tailStart.Find.Execute()

if ($foundTail) {
    $startPos = $tailStart.Start

    # The paragraph that owns this sentence - its end (just before the
    # paragraph mark) is where the sentence finally stops.
    $ownerPara = Get-ParagraphAt $d $startPos
    $paraTextEnd = $ownerPara.Range.End - 1

    # Deleting/replacing text inside a paragraph coalesces every run from
    # the edited point through to the end of the paragraph, so the only
    # way to keep the trailing runs ("å sjekke ", "antallet...", ".")
    # separated exactly like the source is to blow away everything from
    # the edit point to the paragraph end and then rebuild every piece
    # (touched or not) as its own freshly inserted run.
    $victim = $d.Range($startPos, $paraTextEnd)
    $victim.Delete()

    $insPoint = $d.Range($startPos, $startPos)
    $insPoint.InsertAfter(" for opptelling av registreringer per år")
    $insPoint.Collapse(0)

    $insPoint.InsertAfter(". Det kan også være nødvendig ")
    $insPoint.Collapse(0)

    $insPoint.InsertAfter("å sjekke ")
    $insPoint.Collapse(0)

    $insPoint.InsertAfter("antallet registreringer per år på samme måte som med antall mapper per år")
    $insPoint.Collapse(0)

    $insPoint.InsertAfter(".")
    $insPoint.Collapse(0)
}
